$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Formula = '="' + $val.Replace('"', '""') + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

$ws.Range("D2").Value = "26.907.89"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "1.834.77"
$ws.Range("E3").Value = "  -1.68%  "

Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  +0.09%  "

Set-TextValue "D5" "310.60"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("E6").Value = "  +0.06%  "

Set-TextValue "D7" "0.4622"
$ws.Range("E7").Value = "  -0.94%  "

Set-TextValue "D8" "0.3659"
$ws.Range("E8").Value = "  -2.01%  "

Set-TextValue "D9" "0.07157"
$ws.Range("E9").Value = "  -3.23%  "

Set-TextValue "D10" "0.8812"
$ws.Range("E10").Value = "  -0.85%  "

Set-TextValue "D11" "0.07844"
$ws.Range("E11").Value = "  -1.56%  "

Set-TextValue "D12" "19.62"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").Value = "1.832.50"
$ws.Range("E13").Value = "  -1.46%  "

$ws.Range("E14").Value = "  -1.73%  "

Set-TextValue "D15" "6.372"
$ws.Range("E15").Value = "  -3.54%  "

Set-TextValue "D16" "88.72"
$ws.Range("E16").Value = "  -4.35%  "

$ws.Range("E17").Value = "  +0.04%  "

Set-TextValue "D18" "0.000008775"
$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "26.938.31"
$ws.Range("E20").Value = "  -2.06%  "

Set-TextValue "D21" "14.52"
$ws.Range("E21").Value = "  -2.85%  "

$ws.Range("E22").Value = "  -3.31%  "

Set-TextValue "D23" "10.42"
$ws.Range("E23").Value = "  -1.48%  "

Set-TextValue "D24" "1.975"
$ws.Range("E24").Value = "  +5.10%  "

Set-TextValue "D25" "150.83"
$ws.Range("E25").Value = "  -1.57%  "

Set-TextValue "D26" "18.23"
$ws.Range("E26").Value = "  -1.67%  "

Set-TextValue "D27" "2.005"
$ws.Range("E27").Value = "  -4.07%  "

Set-TextValue "D28" "113.66"

Set-TextValue "D29" "4.941"
$ws.Range("E29").Value = "  -4.40%  "

Set-TextValue "D30" "0.08842"
$ws.Range("E30").Value = "  -0.62%  "

Set-TextValue "D31" "3.106"
$ws.Range("E31").Value = "  +3.05%  "

Set-TextValue "D32" "0.7591"
$ws.Range("E32").Value = "  +0.79%  "

Set-TextValue "D33" "4.472"
$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("E34").Value = "  -2.01%  "

Set-TextValue "D35" "2.662"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("E36").Value = "  +1.21%  "

Set-TextValue "D37" "0.01930"
$ws.Range("E37").Value = "  -2.05%  "

Set-TextValue "D38" "2.932"
$ws.Range("E38").Value = "  -1.82%  "

Set-TextValue "D39" "0.05133"
$ws.Range("E39").Value = "  -2.89%  "

Set-TextValue "D40" "6.953"
$ws.Range("E40").Value = "  -3.28%  "

Set-TextValue "D41" "0.4983"
$ws.Range("E41").Value = "  -4.48%  "

Set-TextValue "D42" "0.1598"
$ws.Range("E42").Value = "  -2.88%  "

Set-TextValue "D43" "8.370"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("E44").Value = "  -4.56%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "10.24"
$ws.Range("E45").Value = "  -1.28%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D46" "1.005"
$ws.Range("E46").Value = "  +0.01%  "

Set-TextValue "D47" "102.64"
$ws.Range("E47").Value = "  -0.95%  "

Set-TextValue "D48" "1.614"
$ws.Range("E48").Value = "  -3.24%  "

Set-TextValue "D49" "0.06096"
$ws.Range("E49").Value = "  -2.62%  "

Set-TextValue "D50" "64.62"
$ws.Range("E50").Value = "  -2.00%  "

Set-TextValue "D51" "36.43"
$ws.Range("E51").Value = "  -2.19%  "
